# Rename the inline logo pictures in the document's headers/footers.
#
#   - Pearson logo (footers, both "first page" and "default" footer):
#       image1.png  ->  image2.png
#   - BTec logo (headers, both "first page" and "default" header):
#       image2.jpg  ->  image1.jpg
#
# WdHeaderFooterIndex constants: 1 = wdHeaderFooterPrimary (default),
# 2 = wdHeaderFooterFirstPage.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footers: Pearson logo, image1.png -> image2.png ---
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image2.png"
            }
        }
    }
}

# --- Headers: BTec logo, image2.jpg -> image1.jpg ---
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image1.jpg"
            }
        }
    }
}
